$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 283.33334
$ws.Cells.Item(12, 9).Value = 275
$ws.Cells.Item(12, 10).Value = 300
$ws.Cells.Item(12, 11).Value = 275
$ws.Cells.Item(12, 12).Value = 300
$ws.Cells.Item(12, 13).Value = -105
$ws.Cells.Item(12, 14).Value = -640

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 146.21213
$ws.Cells.Item(15, 9).Value = 146.21213
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 438.63639
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = -269.63639

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 298.4737
$ws.Cells.Item(92, 9).Value = 188.07143
$ws.Cells.Item(92, 10).Value = 607.6
$ws.Cells.Item(92, 11).Value = 188.07143
$ws.Cells.Item(92, 12).Value = 607.6
$ws.Cells.Item(92, 13).Value = 1059.92857

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1355.9412
$ws.Cells.Item(98, 9).Value = 1323.5333
$ws.Cells.Item(98, 10).Value = 1599
$ws.Cells.Item(98, 11).Value = 1323.5333
$ws.Cells.Item(98, 12).Value = 1599
$ws.Cells.Item(98, 13).Value = 174.4666999999999
$ws.Cells.Item(98, 14).Value = -4595

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 1466.125
$ws.Cells.Item(111, 9).Value = 1089
$ws.Cells.Item(111, 10).Value = 2597.5
$ws.Cells.Item(111, 11).Value = 3267
$ws.Cells.Item(111, 12).Value = 7792.5
$ws.Cells.Item(111, 13).Value = -200

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 1355.9412
$ws.Cells.Item(122, 9).Value = 1323.5333
$ws.Cells.Item(122, 10).Value = 1599
$ws.Cells.Item(122, 11).Value = 3970.5999
$ws.Cells.Item(122, 12).Value = 4797
$ws.Cells.Item(122, 13).Value = -1520.5999
$ws.Cells.Item(122, 14).Value = -9697

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 14).Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1667.4166
$ws.Cells.Item(2, 9).Value = 901.55554
$ws.Cells.Item(2, 10).Value = 3965
$ws.Cells.Item(2, 11).Value = 901.55554
$ws.Cells.Item(2, 12).Value = 3965
$ws.Cells.Item(2, 13).Value = -788.55554

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 1868.4849
$ws.Cells.Item(110, 9).Value = 1155.238
$ws.Cells.Item(110, 10).Value = 3116.6667
$ws.Cells.Item(110, 11).Value = 1155.238
$ws.Cells.Item(110, 12).Value = 3116.6667
$ws.Cells.Item(110, 13).Value = 889.7619999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 1667.4166
$ws.Cells.Item(116, 9).Value = 901.55554
$ws.Cells.Item(116, 10).Value = 3965
$ws.Cells.Item(116, 11).Value = 901.55554
$ws.Cells.Item(116, 12).Value = 3965
$ws.Cells.Item(116, 13).Value = 1392.44446

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1667.4166
$ws.Cells.Item(3, 9).Value = 901.55554
$ws.Cells.Item(3, 10).Value = 3965
$ws.Cells.Item(3, 11).Value = 901.55554
$ws.Cells.Item(3, 12).Value = 3965
$ws.Cells.Item(3, 13).Value = -787.55554

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2881.7334
$ws.Cells.Item(105, 9).Value = 2194
$ws.Cells.Item(105, 10).Value = 5632.6665
$ws.Cells.Item(105, 11).Value = 2194
$ws.Cells.Item(105, 12).Value = 5632.6665
$ws.Cells.Item(105, 13).Value = -447

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1859.4286
$ws.Cells.Item(107, 9).Value = 1669.4166
$ws.Cells.Item(107, 10).Value = 2999.5
$ws.Cells.Item(107, 11).Value = 1669.4166
$ws.Cells.Item(107, 12).Value = 2999.5
$ws.Cells.Item(107, 13).Value = 250.5834
$ws.Cells.Item(107, 14).Value = -6839.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(108, 8).Value = 60000
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 60000
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 12).Value = 60000
$ws.Cells.Item(108, 14).Value = -67680

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(130, 8).Value = 55889.9
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 55889.9
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 55889.9
$ws.Cells.Item(130, 14).Value = -65929.89999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 56593.332
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 56593.332
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 56593.332
$ws.Cells.Item(140, 14).Value = -66953.33199999999
$ws.Cells.Item(140, 13).Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 2518.6
$ws.Cells.Item(94, 9).Value = 1949.5
$ws.Cells.Item(94, 10).Value = 2660.875
$ws.Cells.Item(94, 11).Value = 1949.5
$ws.Cells.Item(94, 12).Value = 2660.875
$ws.Cells.Item(94, 13).Value = -1498.5
$ws.Cells.Item(94, 14).Value = -3562.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 979.8889
$ws.Cells.Item(105, 9).Value = 1008.625
$ws.Cells.Item(105, 10).Value = 750
$ws.Cells.Item(105, 11).Value = 1008.625
$ws.Cells.Item(105, 12).Value = 750
$ws.Cells.Item(105, 13).Value = 738.375
$ws.Cells.Item(105, 14).Value = -4244

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2965.1177
$ws.Cells.Item(134, 9).Value = 2402.9092
$ws.Cells.Item(134, 10).Value = 3995.8333
$ws.Cells.Item(134, 11).Value = 7208.7276
$ws.Cells.Item(134, 12).Value = 11987.4999
$ws.Cells.Item(134, 13).Value = -4673.7276
$ws.Cells.Item(134, 14).Value = -17057.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 3237
$ws.Cells.Item(69, 9).Value = 299.5
$ws.Cells.Item(69, 10).Value = 4216.1665
$ws.Cells.Item(69, 11).Value = 898.5
$ws.Cells.Item(69, 12).Value = 12648.4995
$ws.Cells.Item(69, 13).Value = -87.5
$ws.Cells.Item(69, 14).Value = -14270.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(72, 8).Value = 3237
$ws.Cells.Item(72, 9).Value = 299.5
$ws.Cells.Item(72, 10).Value = 4216.1665
$ws.Cells.Item(72, 11).Value = 2695.5
$ws.Cells.Item(72, 12).Value = 37945.4985
$ws.Cells.Item(72, 13).Value = 1360.5
$ws.Cells.Item(72, 14).Value = -46057.4985

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 1061.75
$ws.Cells.Item(109, 9).Value = 785
$ws.Cells.Item(109, 10).Value = 2999
$ws.Cells.Item(109, 11).Value = 2355
$ws.Cells.Item(109, 12).Value = 8997
$ws.Cells.Item(109, 13).Value = -1315

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 13).Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1364.4286
$ws.Cells.Item(102, 9).Value = 1373.909
$ws.Cells.Item(102, 10).Value = 1329.6666
$ws.Cells.Item(102, 11).Value = 1373.909
$ws.Cells.Item(102, 12).Value = 1329.6666
$ws.Cells.Item(102, 13).Value = 248.0909999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 74339.36
$ws.Cells.Item(122, 9).Value = 2198.3
$ws.Cells.Item(122, 10).Value = 254692
$ws.Cells.Item(122, 11).Value = 6594.900000000001
$ws.Cells.Item(122, 12).Value = 764076
$ws.Cells.Item(122, 13).Value = -4144.900000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 9).Value = 6000
$ws.Cells.Item(126, 10).Value = 6333.3335
$ws.Cells.Item(126, 11).Value = 18000
$ws.Cells.Item(126, 12).Value = 19000.0005
$ws.Cells.Item(126, 13).Value = -15530
$ws.Cells.Item(126, 14).Value = -23940.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3126.9
$ws.Cells.Item(132, 9).Value = 3133.5925
$ws.Cells.Item(132, 10).Value = 3066.6667
$ws.Cells.Item(132, 11).Value = 9400.7775
$ws.Cells.Item(132, 12).Value = 9200.000100000001
$ws.Cells.Item(132, 13).Value = -6870.7775

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1500
$ws.Cells.Item(7, 9).Value = 1500
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 1500
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = -1388

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(18, 8).Value = 55333
$ws.Cells.Item(18, 9).Value = 47999.5
$ws.Cells.Item(18, 10).Value = 70000
$ws.Cells.Item(18, 11).Value = 47999.5
$ws.Cells.Item(18, 12).Value = 70000
$ws.Cells.Item(18, 13).Value = -47827.5
$ws.Cells.Item(18, 14).Value = -70344

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(29, 8).Value = 40000
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 40000
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 40000
$ws.Cells.Item(29, 14).Value = -40590

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3408.963
$ws.Cells.Item(40, 9).Value = 3386.2693
$ws.Cells.Item(40, 10).Value = 3999
$ws.Cells.Item(40, 11).Value = 3386.2693
$ws.Cells.Item(40, 12).Value = 3999
$ws.Cells.Item(40, 13).Value = -3250.2693

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 3264.6667
$ws.Cells.Item(68, 9).Value = 3264.6667
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 3264.6667
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = -2515.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 3264.6667
$ws.Cells.Item(71, 9).Value = 3264.6667
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 16323.3335
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = -12579.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 5489
$ws.Cells.Item(100, 9).Value = 5238
$ws.Cells.Item(100, 10).Value = 7999
$ws.Cells.Item(100, 11).Value = 5238
$ws.Cells.Item(100, 12).Value = 7999
$ws.Cells.Item(100, 13).Value = -4697
$ws.Cells.Item(100, 14).Value = -9081

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 6308.25
$ws.Cells.Item(122, 9).Value = 4705.5
$ws.Cells.Item(122, 10).Value = 6995.143
$ws.Cells.Item(122, 11).Value = 14116.5
$ws.Cells.Item(122, 12).Value = 20985.429
$ws.Cells.Item(122, 13).Value = -11666.5
$ws.Cells.Item(122, 14).Value = -25885.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 1500
$ws.Cells.Item(126, 9).Value = 1500
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 4500
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -2030

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 21748.666
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 21748.666
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 21748.666
$ws.Cells.Item(45, 14).Value = -22730.666
$ws.Cells.Item(45, 13).Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1402.3572
$ws.Cells.Item(100, 9).Value = 1387.4445
$ws.Cells.Item(100, 10).Value = 1429.2
$ws.Cells.Item(100, 11).Value = 2774.889
$ws.Cells.Item(100, 12).Value = 2858.4
$ws.Cells.Item(100, 13).Value = -2233.889

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4500
$ws.Cells.Item(132, 9).Value = 4000
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 12000
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -9470
$ws.Cells.Item(132, 14).Value = -20060
